$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NIG(1.0028008308424856, 0.7637242239583764, 0.7988413655512033, 3.254324199005086)"
$ws.Range("C2").Value = "JSU(-1.7029846029909228, 1.311318935784597, 0.7922755769002008, 4.595416393169003)"
$ws.Range("D2").Value = "NIG(0.5212442334468417, 0.3961643278268716, 1.9343484911620235, 2.371216361870894)"
$ws.Range("E2").Value = "NIG(1.6901160107927429, 1.2138031937497542, 4.13471465823856, 6.734669880619813)"
